$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1790.5
$ws.Range("I41").Value = 2387.7144
$ws.Range("J41").Value = 397
$ws.Range("K41").Value = 2387.7144
$ws.Range("L41").Value = 397
$ws.Range("M41").Value = -1947.7144
$ws.Range("N41").Value = -1277

$ws.Range("H42").Value = 148.77777
$ws.Range("I42").Value = 40.8
$ws.Range("J42").Value = 283.75
$ws.Range("K42").Value = 122.4
$ws.Range("L42").Value = 851.25
$ws.Range("M42").Value = 107.6
$ws.Range("N42").Value = -1311.25

$ws.Range("H55").Value = 518.8333
$ws.Range("I55").Value = 592.6
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 592.6
$ws.Range("L55").Value = 150
$ws.Range("M55").Value = -378.6
$ws.Range("N55").Value = -578

$ws.Range("H76").Value = 10000.5
$ws.Range("I76").Value = 10000.5
$ws.Range("K76").Value = 10000.5
$ws.Range("M76").Value = -9685.5

$ws.Range("H79").Value = 10000.5
$ws.Range("I79").Value = 10000.5
$ws.Range("K79").Value = 10000.5
$ws.Range("M79").Value = -8908.5

$ws.Range("H92").Value = 6233.15
$ws.Range("I92").Value = 6831.6665
$ws.Range("K92").Value = 6831.6665
$ws.Range("M92").Value = -5583.6665

$ws.Range("H106").Value = 1395.7
$ws.Range("I106").Value = 1108.5
$ws.Range("K106").Value = 1108.5
$ws.Range("M106").Value = -477.5

$ws.Range("H107").Value = 419.42856
$ws.Range("I107").Value = 387.3793
$ws.Range("K107").Value = 387.3793
$ws.Range("M107").Value = 1532.6207

$ws.Range("H111").Value = 13260.375
$ws.Range("I111").Value = 13260.375
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 39781.125
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -36714.125
$ws.Range("N111").ClearContents()

$ws.Range("H116").Value = 10056
$ws.Range("J116").Value = 8640.909
$ws.Range("L116").Value = 8640.909
$ws.Range("N116").Value = -15524.909

$ws.Range("H132").Value = 30395990
$ws.Range("I132").Value = 37039220
$ws.Range("J132").Value = 501450
$ws.Range("K132").Value = 111117660
$ws.Range("L132").Value = 1504350
$ws.Range("M132").Value = -111115130
$ws.Range("N132").Value = -1509410

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5366.3
$ws.Range("I61").Value = 4915.857
$ws.Range("J61").Value = 6417.3335
$ws.Range("K61").Value = 4915.857
$ws.Range("L61").Value = 6417.3335
$ws.Range("M61").Value = -4703.857
$ws.Range("N61").Value = -6841.3335

$ws.Range("H74").Value = 6747.0312
$ws.Range("I74").Value = 996.2727
$ws.Range("J74").Value = 19398.7
$ws.Range("K74").Value = 996.2727
$ws.Range("L74").Value = 19398.7
$ws.Range("M74").Value = -122.2727
$ws.Range("N74").Value = -21146.7

$ws.Range("H77").Value = 6747.0312
$ws.Range("I77").Value = 996.2727
$ws.Range("J77").Value = 19398.7
$ws.Range("K77").Value = 4981.363499999999
$ws.Range("L77").Value = 96993.5
$ws.Range("M77").Value = -613.3634999999995
$ws.Range("N77").Value = -105729.5

$ws.Range("H102").Value = 3642.5715
$ws.Range("I102").Value = 3499.6843
$ws.Range("K102").Value = 3499.6843
$ws.Range("M102").Value = -1877.6843

$ws.Range("H122").Value = 1680.4
$ws.Range("I122").Value = 1514.7142
$ws.Range("K122").Value = 4544.142599999999
$ws.Range("M122").Value = -2094.142599999999

$ws.Range("H136").Value = 5366.3
$ws.Range("I136").Value = 4915.857
$ws.Range("J136").Value = 6417.3335
$ws.Range("K136").Value = 14747.571
$ws.Range("L136").Value = 19252.0005
$ws.Range("M136").Value = -12197.571
$ws.Range("N136").Value = -24352.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6118
$ws.Range("I105").Value = 8296.666999999999
$ws.Range("J105").Value = 2850
$ws.Range("K105").Value = 8296.666999999999
$ws.Range("L105").Value = 2850
$ws.Range("M105").Value = -6549.666999999999
$ws.Range("N105").Value = -6344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38051
$ws.Range("I31").Value = 49024.523
$ws.Range("J31").Value = 9245.5
$ws.Range("K31").Value = 49024.523
$ws.Range("L31").Value = 9245.5
$ws.Range("M31").Value = -48729.523
$ws.Range("N31").Value = -9835.5

$ws.Range("H34").Value = 38051
$ws.Range("I34").Value = 49024.523
$ws.Range("J34").Value = 9245.5
$ws.Range("K34").Value = 49024.523
$ws.Range("L34").Value = 9245.5
$ws.Range("M34").Value = -48822.523
$ws.Range("N34").Value = -9649.5

$ws.Range("H132").Value = 4346.0415
$ws.Range("I132").Value = 4305.048
$ws.Range("K132").Value = 12915.144
$ws.Range("M132").Value = -10385.144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 9025995
$ws.Range("I9").Value = 19800800
$ws.Range("J9").Value = 46991.668
$ws.Range("K9").Value = 59402400
$ws.Range("L9").Value = 140975.004
$ws.Range("M9").Value = -59402176
$ws.Range("N9").Value = -141423.004

$ws.Range("H11").Value = 370
$ws.Range("J11").Value = 212.5
$ws.Range("L11").Value = 637.5
$ws.Range("N11").Value = -917.5

$ws.Range("H64").Value = 1000
$ws.Range("J64").Value = 1000
$ws.Range("L64").Value = 3000
$ws.Range("N64").Value = -3540

$ws.Range("H67").Value = 1000
$ws.Range("J67").Value = 1000
$ws.Range("L67").Value = 3000
$ws.Range("N67").Value = -4872

$ws.Range("H69").Value = 3801.5
$ws.Range("I69").Value = 3801.5
$ws.Range("K69").Value = 11404.5
$ws.Range("M69").Value = -10593.5

$ws.Range("H72").Value = 3801.5
$ws.Range("I72").Value = 3801.5
$ws.Range("K72").Value = 34213.5
$ws.Range("M72").Value = -30157.5

$ws.Range("H98").Value = 539
$ws.Range("I98").Value = 559
$ws.Range("J98").Value = 499
$ws.Range("K98").Value = 1677
$ws.Range("L98").Value = 1497
$ws.Range("M98").Value = -179
$ws.Range("N98").Value = -4493

$ws.Range("H111").Value = 3655.3333
$ws.Range("I111").Value = 2379.6
$ws.Range("J111").Value = 5250
$ws.Range("K111").Value = 7138.799999999999
$ws.Range("L111").Value = 15750
$ws.Range("M111").Value = -4071.799999999999
$ws.Range("N111").Value = -21884

$ws.Range("H118").Value = 3412.4167
$ws.Range("I118").Value = 999.5
$ws.Range("J118").Value = 3895
$ws.Range("K118").Value = 2998.5
$ws.Range("L118").Value = 11685
$ws.Range("M118").Value = -1755.5
$ws.Range("N118").Value = -14171

$ws.Range("H129").Value = 936.5
$ws.Range("I129").Value = 748.6667
$ws.Range("K129").Value = 2246.0001
$ws.Range("M129").Value = 2753.9999

$ws.Range("H131").Value = 47892.273
$ws.Range("I131").Value = 84210.75
$ws.Range("J131").Value = 4310.1
$ws.Range("K131").Value = 252632.25
$ws.Range("L131").Value = 12930.3
$ws.Range("M131").Value = -247592.25
$ws.Range("N131").Value = -23010.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 39933.332
$ws.Range("J123").Value = 39933.332
$ws.Range("L123").Value = 39933.332
$ws.Range("N123").Value = -44833.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1196.6666
$ws.Range("J22").Value = 1196.6666
$ws.Range("L22").Value = 1196.6666
$ws.Range("N22").Value = -1786.6666

$ws.Range("H27").Value = 1196.6666
$ws.Range("J27").Value = 1196.6666
$ws.Range("L27").Value = 1196.6666
$ws.Range("N27").Value = -1410.6666

$ws.Range("H74").Value = 70000
$ws.Range("J74").Value = 70000
$ws.Range("L74").Value = 70000
$ws.Range("N74").Value = -71996

$ws.Range("H77").Value = 70000
$ws.Range("J77").Value = 70000
$ws.Range("L77").Value = 210000
$ws.Range("N77").Value = -219984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 670.3684
$ws.Range("I14").Value = 646.35297
$ws.Range("K14").Value = 646.35297
$ws.Range("M14").Value = -478.35297

$ws.Range("H100").Value = 1496.2106
$ws.Range("I100").Value = 984
$ws.Range("J100").Value = 1795
$ws.Range("K100").Value = 1968
$ws.Range("L100").Value = 3590
$ws.Range("M100").Value = -1427
$ws.Range("N100").Value = -4672

$ws.Range("H113").Value = 734.9394
$ws.Range("I113").Value = 666.9
$ws.Range("J113").Value = 839.61536
$ws.Range("K113").Value = 2000.7
$ws.Range("L113").Value = 2518.84608
$ws.Range("M113").Value = 169.3000000000002
$ws.Range("N113").Value = -6858.84608

$ws.Range("H122").Value = 3638.2
$ws.Range("I122").Value = 3469.5
$ws.Range("K122").Value = 10408.5
$ws.Range("M122").Value = -7958.5

$ws.Range("H136").Value = 2465.65
$ws.Range("J136").Value = 2053.75
$ws.Range("L136").Value = 6161.25
$ws.Range("N136").Value = -11261.25
